# "E suite.xlsx" - update the Results column (D) for rows 9-11 on the
# "Test Cases" sheet from "FAIL" to "SKIP".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D9").Value = "SKIP"
$ws.Range("D10").Value = "SKIP"
$ws.Range("D11").Value = "SKIP"
